# fix: more cleaning for excels, countries
#
# Row 2 (Bowen MD, Peters CJ, Nichol ST., 1997) has a "Country" entry
# (column Q) that incorrectly uses a slash to join two countries:
# "Mozambique/Zimbabwe". Clean it up to use a comma like the other
# country lists in the sheet ("Mozambique, Zimbabwe").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = "West Africa, Central African Republic, Mozambique, Zimbabwe, South America"

# Leave the cursor/selection where the author ended up after this edit.
$ws.Activate()
$ws.Range("O3").Select()

$wb.Save()
